$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('AK_1.png', 0.539, 0.002, 0.982, 0.333, 'Tidak Diketahui', 'Salah'),
    @('AK_2.png', 0.896, 0.003, 0.988, 0.4, 'Tidak Diketahui', 'Salah'),
    @('AK_3.png', 0.86, 0.003, 0.979, 0.8, 'Akhlak Kamiswara', 'Benar'),
    @('AK_4.png', 0.538, 0.002, 0.976, 0.533, 'Akhlak Kamiswara', 'Benar'),
    @('AK_5.png', 0.537, 0.002, 0.99, 0.333, 'Tidak Diketahui', 'Salah'),
    @('MIB_1.png', 1.613, 0.005, 0.976, 0.667, 'Muhammad Iqbal Baqi', 'Benar'),
    @('MIB_2.png', 1.214, 0.004, 0.973, 0.467, 'Tidak Diketahui', 'Salah'),
    @('MIB_3.png', 1.765, 0.006, 0.999, 1, 'Muhammad Iqbal Baqi', 'Benar'),
    @('MIB_4.png', 1.178, 0.004, 0.985, 0.733, 'Muhammad Iqbal Baqi', 'Benar'),
    @('MIB_5.png', 1.587, 0.005, 0.983, 0.6, 'Muhammad Iqbal Baqi', 'Benar'),
    @('AAH_1.png', 0.876, 0.003, 0.975, 0.533, 'Andrea Ayunove Hutami', 'Benar'),
    @('AAH_2.png', 1.104, 0.004, 0.992, 0.6, 'Andrea Ayunove Hutami', 'Benar'),
    @('AAH_3.png', 0.918, 0.003, 0.973, 0.667, 'Andrea Ayunove Hutami', 'Benar'),
    @('TI_1.png', 0.947, 0.003, 0.949, 0.533, 'Toni Ismail', 'Benar'),
    @('TI_2.png', 0.953, 0.003, 0.989, 0.333, 'Tidak Diketahui', 'Salah'),
    @('TI_3.png', 0.757, 0.003, 0.991, 0.6, 'Toni Ismail', 'Benar'),
    @('TI_4.png', 0.751, 0.003, 0.974, 0.733, 'Toni Ismail', 'Benar'),
    @('TI_5.png', 1.069, 0.004, 0.992, 0.467, 'Tidak Diketahui', 'Salah'),
    @('RAS_1.png', 0.6830000000000001, 0.002, 0.975, 0.333, 'Tidak Diketahui', 'Salah'),
    @('RAS_2.png', 1.08, 0.004, 0.966, 0.6, 'Ridha Ayu Salsabila', 'Benar'),
    @('RAS_3.png', 0.597, 0.002, 0.982, 0.4, 'Tidak Diketahui', 'Salah'),
    @('RAS_4.png', 1.315, 0.004, 0.975, 0.267, 'Tidak Diketahui', 'Salah'),
    @('RAS_5.png', 1.188, 0.004, 0.971, 0.8, 'Ridha Ayu Salsabila', 'Benar'),
    @('RR_1.png', 1.54, 0.005, 0.983, 0.333, 'Tidak Diketahui', 'Salah'),
    @('RR_2.png', 1.362, 0.004, 0.987, 0.333, 'Tidak Diketahui', 'Salah'),
    @('RR_3.png', 1.155, 0.004, 0.906, 0.4, 'Tidak Diketahui', 'Salah'),
    @('RR_4.png', 1.417, 0.005, 0.989, 0.4, 'Tidak Diketahui', 'Salah'),
    @('RR_5.png', 1.483, 0.005, 0.983, 0.667, 'Rafiqo Rapitasari', 'Benar'),
    @('AR_1.png', 0.886, 0.003, 0.965, 0.733, 'Arizli Romadhon', 'Benar'),
    @('GA_1.png', 1.714, 0.006, 0.983, 0.8, 'Gege Ardiyansyah', 'Benar'),
    @('GA_2.png', 0.6879999999999999, 0.002, 0.99, 0.267, 'Tidak Diketahui', 'Salah'),
    @('GA_3.png', 0.743, 0.002, 0.99, 0.333, 'Tidak Diketahui', 'Salah'),
    @('FY_1.png', 1.141, 0.004, 0.969, 0.333, 'Tidak Diketahui', 'Salah'),
    @('FY_2.png', 1.693, 0.006, 0.98, 0.333, 'Tidak Diketahui', 'Salah'),
    @('FY_3.png', 1.516, 0.005, 0.991, 0.4, 'Tidak Diketahui', 'Salah'),
    @('FY_4.png', 1.262, 0.004, 0.984, 0.267, 'Tidak Diketahui', 'Salah'),
    @('TO_1.png', 0.819, 0.003, 0.971, 0.4, 'Tidak Diketahui', 'Salah'),
    @('TO_2.png', 1.066, 0.004, 0.986, 0.533, 'Tiara Oktavian', 'Benar'),
    @('TO_3.png', 0.841, 0.003, 0.983, 0.6, 'Tiara Oktavian', 'Benar'),
    @('TO_4.png', 7.237, 0.024, 0.9320000000000001, 0.333, 'Tidak Diketahui', 'Salah'),
    @('TO_5.png', 7.251, 0.024, 0.9370000000000001, 0.267, 'Tidak Diketahui', 'Salah'),
    @('TD_1.png', 3.814, 0.013, 0.894, 0.467, 'Tidak Diketahui', 'Benar'),
    @('TD_2.png', 5.068, 0.017, 0.885, 0.267, 'Tidak Diketahui', 'Benar'),
    @('TD_3.png', 1.071, 0.004, 0.955, 0.267, 'Tidak Diketahui', 'Benar'),
    @('TD_4.png', 1.009, 0.003, 0.987, 0.267, 'Tidak Diketahui', 'Benar')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
